$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix KeyError in portfolio mapping column name handling:
# rename the column header in B1 from "level1accountname" to "Fund_Code"
$ws.Range("B1").Value = "Fund_Code"

# Update the active selection to match where the user ended up (I9)
$ws.Range("I9").Select()
